$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new forecast vintage (column BB, the 54th column) is added to the staircase.
# Column BB1 holds the new vintage's date header; BB2..BB83 hold that vintage's
# QoQ forecast values. For the rows where the series had already flattened out
# (rows 2-70) the new vintage simply repeats the prior (BA) value unchanged, so
# those are populated by copying the BA cell (value + number format/border/bold
# style) one column to the right. Rows 71-82 get genuinely new forecast values,
# and a brand-new row 83 (next quarter, date 46934) is appended with only the
# BB83 forecast populated (mirrors how each new vintage column starts business
# further out than the last and the staircase only back-fills so far).

# --- BB1: new vintage date header. Copy BA1's style (bold/border/date format) then overwrite with the new date serial.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# --- BB2:BB70 repeat the BA value verbatim -> copy BA -> BB cell by cell (keeps value identical, no style to carry since these are plain data cells).
for ($r = 2; $r -le 70; $r++) {
    $ws.Range("BA$r").Copy($ws.Range("BB$r"))
}

# --- BB71:BB82: new forecast values for this vintage (differ from the BA column).
$ws.Range("BB71").Value = 0.5597354586130052
$ws.Range("BB72").Value = 0.1368731201391853
$ws.Range("BB73").Value = -0.2551464291630765
$ws.Range("BB74").Value = 0.1643375991815219
$ws.Range("BB75").Value = 0.1643375991815219
$ws.Range("BB76").Value = 0.1643375991815219
$ws.Range("BB77").Value = 0.1643375991815219
$ws.Range("BB78").Value = 0.1643375991815219
$ws.Range("BB79").Value = 0.1643375991815219
$ws.Range("BB80").Value = 0.1643375991815219
$ws.Range("BB81").Value = 0.1643375991815219
$ws.Range("BB82").Value = 0.1643375991815219

# --- New row 83 (next quarter's date, 46934). Column A gets the same date style as the rest of column A; only BB83 (this new vintage) has a forecast for that quarter.
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.1643375991815219

Write-Host "edit complete"
